$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current last data row (row 7), pushing the
# existing "Desc" row down to row 9. Excel inherits formatting from the row
# above on insert, matching the style (s="1") already used by the sibling
# data rows.
$ws.Rows("7:8").Insert()

# New row 7: "Icon" field definition
$ws.Cells.Item(7, 1).Value = "Icon"
$ws.Cells.Item(7, 2).Value = "string"
$ws.Cells.Item(7, 3).Value = $false
$ws.Cells.Item(7, 4).Value = $false
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = $false
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = "Friend"
$ws.Cells.Item(7, 10).Value = "图标"

# New row 8: "ShowName" field definition
$ws.Cells.Item(8, 1).Value = "ShowName"
$ws.Cells.Item(8, 2).Value = "string"
$ws.Cells.Item(8, 3).Value = $false
$ws.Cells.Item(8, 4).Value = $false
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = $false
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = "Friend"
$ws.Cells.Item(8, 10).Value = "名字"

# Move the cursor/selection to reflect the author's final position in the sheet.
$ws.Range("C13").Select() | Out-Null
